$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column C for rows 2-12
# from 45170 (2023-09-01) to 45174 (2023-09-05), keeping existing date formatting.
$newDate = [DateTime]::FromOADate(45174)

for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
